$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Michele Bertolini"
$ws.Range("B5").Value = "Nicola Togni | RSA United"
$ws.Range("C5").Value = "Alessandro Fait | RSA United"
$ws.Range("D5").Value = "Leonardo  Parisi  | MediaserT"
$ws.Range("E5").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("F5").Value = "Luca Lasta | La Contea FC"
